$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Skill tree / tech text: repurpose the "Image Gen ID" / "Image GEN ID" column (K)
# into a "TECH12" marker, and give row 4's previously-numeric placeholder a
# "NULL" text marker instead of 0.
$ws.Range("K2").Value = "TECH12"
$ws.Range("K3").Value = "TECH12"
$ws.Range("K4").Value = "NULL"

# Move the active selection to K4, matching where the edits were made.
$null = $ws.Range("K4").Select()
